$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D numeric-looking strings remain plain text (matches source formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.298.09"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.400.93"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.62"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.85"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +7.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.585"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.38"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "680.58"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.948.50"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.412.13"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.403.65"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.70"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.30"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.04"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.88"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.89"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.46"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("E31").Value = "  +11.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "556.62"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.02"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.606.90"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.23"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0745"
$ws.Range("E40").Value = "  +10.13%  "
$ws.Range("E41").Value = "  +3.69%  "
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.67"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.35"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.65"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("E51").Value = "  -0.92%  "
